$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) survey_3: remove the stray "0" value that was sitting in B13
#    (row 13 / column B) next to the "insq26a" variable name.
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("survey_3")
$ws4.Range("B13").ClearContents()

# ------------------------------------------------------------------
# 2) Duplicate remove_vars (in its current, pre-"completedwaves" state)
#    into a new sheet "remove_vars_contact" that keeps only the first
#    120 rows and drops the carried-over sort state.
# ------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("remove_vars")
$ws6.Copy($null, $ws6)
$ws7 = $wb.Worksheets.Item("remove_vars (2)")
$ws7.Name = "remove_vars_contact"
$null = $ws7.Range("A121:A161").EntireRow.Delete()

# ------------------------------------------------------------------
# 3) remove_vars: append a new row for the "completedwaves" variable
#    (it already exists at the top of the sheet as row 2 - duplicate
#    entry added at the bottom, matching the style of the rows above
#    it).
# ------------------------------------------------------------------
$ws6.Range("A162").Value = "completedwaves"
$ws6.Range("A161").Copy()
$null = $ws6.Range("A162").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Restore / update the on-screen selection state for each sheet
#    that moved.
# ------------------------------------------------------------------
$null = $ws7.Columns("A").Select()

$null = $ws6.Range("A23").Select()

$ws4.Activate()
$null = $ws4.Range("B13").Select()
